$wb = $excel.ActiveWorkbook

# hunk 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 105.111115
$ws.Range("I5").Value = 110.15385
$ws.Range("J5").Value = 92
$ws.Range("K5").Value = 110.15385
$ws.Range("L5").Value = 92
$ws.Range("M5").Value = 4.846149999999994
$ws.Range("N5").Value = -322

# hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 179.3
$ws.Range("I33").Value = 115.69231
$ws.Range("J33").Value = 297.42856
$ws.Range("K33").Value = 115.69231
$ws.Range("L33").Value = 297.42856
$ws.Range("M33").Value = 113.30769
$ws.Range("N33").Value = -755.4285600000001

# hunk 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 143.56522
$ws.Range("I55").Value = 114.57143
$ws.Range("J55").Value = 188.66667
$ws.Range("K55").Value = 114.57143
$ws.Range("L55").Value = 188.66667
$ws.Range("M55").Value = 99.42856999999999
$ws.Range("N55").Value = -616.6666700000001

# hunk 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1322.591
$ws.Range("I125").Value = 465
$ws.Range("J125").Value = 1574.8235
$ws.Range("K125").Value = 4185
$ws.Range("L125").Value = 14173.4115
$ws.Range("M125").Value = -1725
$ws.Range("N125").Value = -19093.4115

# hunk 4: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384

# hunk 5: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 110.5
$ws.Range("I5").Value = 113.333336
$ws.Range("J5").Value = 102
$ws.Range("K5").Value = 113.333336
$ws.Range("L5").Value = 102
$ws.Range("M5").Value = -1.333336000000003
$ws.Range("N5").Value = -326

# hunk 6: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 800
$ws.Range("I26").Value = 800
$ws.Range("K26").Value = 800
$ws.Range("M26").Value = -470

# hunk 7: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 37040524
$ws.Range("I132").Value = 58825296
$ws.Range("J132").Value = 6414.2
$ws.Range("K132").Value = 176475888
$ws.Range("L132").Value = 19242.6
$ws.Range("M132").Value = -176473358
$ws.Range("N132").Value = -24302.6

# hunk 8: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 110.5
$ws.Range("I4").Value = 113.333336
$ws.Range("J4").Value = 102
$ws.Range("K4").Value = 113.333336
$ws.Range("L4").Value = 102
$ws.Range("M4").Value = 1.666663999999997
$ws.Range("N4").Value = -332

# hunk 9: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 19260
$ws.Range("J81").Value = 19260
$ws.Range("L81").Value = 19260
$ws.Range("N81").Value = -21382

# hunk 10: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 19260
$ws.Range("J84").Value = 19260
$ws.Range("L84").Value = 57780
$ws.Range("N84").Value = -68388

# hunk 11: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 26.875
$ws.Range("I7").Value = 25.571428
$ws.Range("J7").Value = 36
$ws.Range("K7").Value = 25.571428
$ws.Range("L7").Value = 36
$ws.Range("M7").Value = 87.428572
$ws.Range("N7").Value = -262

# hunk 12: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4512.8486
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4512.8486
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4512.8486
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5102.8486

# hunk 13: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4512.8486
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4512.8486
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 4512.8486
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -4916.8486

# hunk 14: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 11159.75
$ws.Range("I35").Value = 2450
$ws.Range("J35").Value = 25676
$ws.Range("K35").Value = 2450
$ws.Range("L35").Value = 25676
$ws.Range("M35").Value = -2156
$ws.Range("N35").Value = -26264

# hunk 15: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 23707.6
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 23707.6
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 23707.6
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -24729.6

# hunk 16: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 32795.832
$ws.Range("J88").Value = 32795.832
$ws.Range("L88").Value = 32795.832
$ws.Range("N88").Value = -33607.832

# hunk 17: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 32795.832
$ws.Range("J91").Value = 32795.832
$ws.Range("L91").Value = 32795.832
$ws.Range("N91").Value = -35603.832

# hunk 18: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1674.1111
$ws.Range("I99").Value = 1553.3334
$ws.Range("J99").Value = 1915.6666
$ws.Range("K99").Value = 1553.3334
$ws.Range("L99").Value = 1915.6666
$ws.Range("M99").Value = -55.33339999999998
$ws.Range("N99").Value = -4911.6666

# hunk 19: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1674.1111
$ws.Range("I126").Value = 1553.3334
$ws.Range("J126").Value = 1915.6666
$ws.Range("K126").Value = 4660.0002
$ws.Range("L126").Value = 5746.9998
$ws.Range("M126").Value = -2190.0002
$ws.Range("N126").Value = -10686.9998

# hunk 20: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 35473.332
$ws.Range("J138").Value = 35473.332
$ws.Range("L138").Value = 35473.332
$ws.Range("N138").Value = -45753.332

# hunk 21: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 54499.5
$ws.Range("J139").Value = 54499.5
$ws.Range("L139").Value = 54499.5
$ws.Range("N139").Value = -64779.5

# hunk 22: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1642888.1
$ws.Range("I2").Value = 83358.25
$ws.Range("J2").Value = 4761948
$ws.Range("K2").Value = 500149.5
$ws.Range("L2").Value = 28571688
$ws.Range("M2").Value = -500036.5
$ws.Range("N2").Value = -28571914

# hunk 23: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 69030.39999999999
$ws.Range("I129").Value = 1653.3334
$ws.Range("J129").Value = 85874.664
$ws.Range("K129").Value = 4960.0002
$ws.Range("L129").Value = 257623.992
$ws.Range("M129").Value = 39.9997999999996
$ws.Range("N129").Value = -267623.992

# hunk 24: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.25
$ws.Range("I2").Value = 25.375
$ws.Range("J2").Value = 47.125
$ws.Range("K2").Value = 25.375
$ws.Range("L2").Value = 47.125
$ws.Range("M2").Value = 87.625
$ws.Range("N2").Value = -273.125

# hunk 25: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2859.5454
$ws.Range("I80").Value = 2730.4
$ws.Range("J80").Value = 2967.1667
$ws.Range("K80").Value = 2730.4
$ws.Range("L80").Value = 2967.1667
$ws.Range("M80").Value = -1732.4
$ws.Range("N80").Value = -4963.1667

# hunk 26: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2859.5454
$ws.Range("I83").Value = 2730.4
$ws.Range("J83").Value = 2967.1667
$ws.Range("K83").Value = 13652
$ws.Range("L83").Value = 14835.8335
$ws.Range("M83").Value = -8660
$ws.Range("N83").Value = -24819.8335

# hunk 27: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 16590.412
$ws.Range("J136").Value = 16590.412
$ws.Range("L136").Value = 49771.236
$ws.Range("N136").Value = -54871.236

# hunk 28: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1141.7858
$ws.Range("I55").Value = 225.25
$ws.Range("J55").Value = 1508.4
$ws.Range("K55").Value = 225.25
$ws.Range("L55").Value = 1508.4
$ws.Range("M55").Value = -52.25
$ws.Range("N55").Value = -1854.4

# hunk 29: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1631.3572
$ws.Range("I68").Value = 1619.9166
$ws.Range("K68").Value = 1619.9166
$ws.Range("M68").Value = -870.9166

# hunk 30: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1631.3572
$ws.Range("I71").Value = 1619.9166
$ws.Range("K71").Value = 8099.583000000001
$ws.Range("M71").Value = -4355.583000000001

# hunk 31: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4476
$ws.Range("I132").Value = 4312.143
$ws.Range("J132").Value = 5049.5
$ws.Range("K132").Value = 12936.429
$ws.Range("L132").Value = 15148.5
$ws.Range("M132").Value = -10406.429
$ws.Range("N132").Value = -20208.5

# hunk 32: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1052.409
$ws.Range("J96").Value = 1105.5
$ws.Range("L96").Value = 1105.5
$ws.Range("N96").Value = -3851.5

# hunk 33: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1381.7971
$ws.Range("I132").Value = 1523.5454
$ws.Range("K132").Value = 4570.6362
$ws.Range("M132").Value = -2040.6362
